$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Logs" sheet: append a new row (row 20) describing the latest e-mail.
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A20").Value = "Aanmelden voor beurs"
$logs.Range("B20").Value = "mailmind.test@zohomail.eu"
$logs.Range("C20").Value = "Wij willen graag deelnemen aan de komende vakbeurs."
$logs.Range("D20").Value = "Samenwerking / Partnerverzoek"
$logs.Range("E20").Value = "Beste,`nBedankt voor uw interesse in deelname aan de komende vakbeurs. Voor meer informatie over de mogelijkheden en beschikbaarheid verzoek ik u vriendelijk om contact op te nemen met onze evenementencoördinator via [contactgegevens]. Hij/zij zal u verder kunnen informeren over de beschikbare standruimte, kosten en eventuele inschrijvingsprocedure.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$logs.Range("F20").Value = "2025-06-22 18:50:16"
$logs.Range("G20").Value = "Ja"

# Undo the automatic row-height ("wrap text") adjustment that the engine
# applies when a multi-line value is written, so the row stays at the
# default height like every other row in the sheet.
$logs.Rows.Item(20).AutoFit()

# The conditional formatting ranges on the Categorie/Beantwoord columns must
# grow by one row to keep covering the whole data range.
$dConditions = $logs.Range("D2:D19").FormatConditions
for ($i = 1; $i -le $dConditions.Count; $i++) {
    $dConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D20"))
}

$gConditions = $logs.Range("G2:G19").FormatConditions
for ($i = 1; $i -le $gConditions.Count; $i++) {
    $gConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G20"))
}

# ---------------------------------------------------------------------------
# 2. "Dashboard" sheet: the per-category counts table is regenerated once the
#    new e-mail is taken into account. "Samenwerking / Partnerverzoek" now
#    has 2 occurrences, so it moves up to row 6, and every category that used
#    to sit below it shifts down by one row.
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B6").Value = 2

$dash.Range("A7").Value = "Klacht / Probleem"
$dash.Range("B7").Value = 1

$dash.Range("A8").Value = "Uitnodiging / Evenement"
$dash.Range("B8").Value = 1

$dash.Range("A9").Value = "Offerte / Prijsaanvraag"
$dash.Range("B9").Value = 1

$dash.Range("A10").Value = "Openingstijden / Locatie"
$dash.Range("B10").Value = 1
